# Add a new "2021" data column (column R) to the 1.5.1 indicator sheet,
# mirroring the formatting of the existing "2020" column (Q) and filling
# in the values for each region / sex breakdown row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 is a blank separator row - just needs R3 to carry the same
# (bottom-border) formatting as the rest of the row.
$ws.Range("Q3").Copy() | Out-Null
$ws.Range("R3").PasteSpecial(-4122) | Out-Null

# Map of row number -> value to place in column R (numbers, or "-" for
# the cells that have no data, matching column Q's placeholder style).
$values = @{
    4  = 2021
    5  = 109
    6  = 74
    7  = 35
    8  = 36
    9  = 35
    10 = 1
    11 = 15
    12 = 8
    13 = 7
    14 = 12
    15 = 7
    16 = 5
    17 = "-"
    18 = "-"
    19 = "-"
    20 = 17
    21 = 8
    22 = 9
    23 = 9
    24 = 7
    25 = 2
    26 = 20
    27 = 9
    28 = 11
    29 = "-"
    30 = "-"
    31 = "-"
    32 = "-"
    33 = "-"
    34 = "-"
}

foreach ($r in 4..34) {
    $srcCell = $ws.Range("Q$r")
    $dstCell = $ws.Range("R$r")

    # Copy Q's formatting onto R so the new column matches the table style.
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null

    $dstCell.Value = $values[$r]
}

# Move the selection off the stale R13 reference now that the column has
# real data in it.
$ws.Range("A1").Select() | Out-Null
